# "version with 2 ref files"
# - Clear out the old "Inandout" text from B2 (keep its red-font style)
# - Add a new reference column E with per-row lookup values
# - Append two new reference rows (B12/B13) with new text values
# - Move the active selection to B14 (just past the new data)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: drop the stale "Inandout" label, keep the cell's (red) style
$ws.Range("B2").ClearContents()

# New column E: per-row reference values
$ws.Range("E2").Value = 5
$ws.Range("E3").Value = 6
$ws.Range("E4").Value = 7
$ws.Range("E5").Value = 3
$ws.Range("E6").Value = 4
$ws.Range("E7").Value = 0
$ws.Range("E8").Value = 1
$ws.Range("E9").Value = 2
$ws.Range("E10").Value = 8
$ws.Range("E11").Value = 9

# Two new reference rows
$ws.Range("B12").Value = "sdcasdc"
$ws.Range("B13").Value = "sdcasddd"

# Leave the selection where the user would land next
$ws.Range("B14").Select()
